$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix duplicated "类类" -> "类" typo in row-1 headers ---
$ws.Range("B1").Value = "其他食品类城市居民消费价格指数(上年=100)"
$ws.Range("D1").Value = "奶类城市居民消费价格指数(上年=100)"
$ws.Range("E1").Value = "干鲜瓜果类城市居民消费价格指数(上年=100)"
$ws.Range("H1").Value = "畜肉类城市居民消费价格指数(上年=100)"
$ws.Range("I1").Value = "禽肉类城市居民消费价格指数(上年=100)"
$ws.Range("K1").Value = "糖果糕点类城市居民消费价格指数(上年=100)"
$ws.Range("N1").Value = "薯类城市居民消费价格指数(上年=100)"
$ws.Range("O1").Value = "蛋类城市居民消费价格指数(上年=100)"
$ws.Range("Q1").Value = "豆类城市居民消费价格指数(上年=100)"

# --- Add the 2021 year row (row 7), matching the style of the A-column year label cells ---
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 100.4
$ws.Range("C7").Value = 102
$ws.Range("D7").Value = 101.9
$ws.Range("E7").Value = 102.1
$ws.Range("F7").Value = 109
$ws.Range("G7").Value = 101.6
$ws.Range("H7").Value = 83.7
$ws.Range("I7").Value = 96.5
$ws.Range("J7").Value = 101.1
$ws.Range("K7").Value = 101.6
$ws.Range("L7").Value = 101.2
$ws.Range("M7").Value = 105
$ws.Range("N7").Value = 99.40000000000001
$ws.Range("O7").Value = 110.2
$ws.Range("P7").Value = 101.2
$ws.Range("Q7").Value = 106.6
$ws.Range("R7").Value = 100
$ws.Range("S7").Value = 99
$ws.Range("T7").Value = 107.4
$ws.Range("U7").Value = 102.9
$ws.Range("V7").Value = 105.6

# --- Add the 2022 year row (row 8); only the 食品烟酒类 figure (R8) has been published so far ---
$ws.Range("A6").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "2022年"
$ws.Range("R8").Value = 102.6
